# Update the NATMI ligand-receptor pairwise output (Col18a1-Ptprs) with
# recalculated TPM-derived values (per commit: "update scripts wuth new tpm").
# Only numeric statistic columns (E:T, excluding K/L which are unchanged)
# for data rows 2-10 are updated to the newly recomputed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5586043333333334
$ws.Range("H2").Value = 1.675813
$ws.Range("I2").Value = 0.01643366487114074
$ws.Range("J2").Value = 0.01643366487114074
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 0.9275183657576666
$ws.Range("R2").Value = 8.347665291818998
$ws.Range("S2").Value = 0.0006104761040176394
$ws.Range("T2").Value = 0.0006104761040176392
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5586043333333334
$ws.Range("H3").Value = 1.675813
$ws.Range("I3").Value = 0.01643366487114074
$ws.Range("J3").Value = 0.01643366487114074
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("Q3").Value = 14.06060751636144
$ws.Range("R3").Value = 126.545467647253
$ws.Range("S3").Value = 0.009254441975062875
$ws.Range("T3").Value = 0.009254441975062875
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5586043333333334
$ws.Range("H4").Value = 1.675813
$ws.Range("I4").Value = 0.01643366487114074
$ws.Range("J4").Value = 0.01643366487114074
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 9.980133947178334
$ws.Range("R4").Value = 89.82120552460501
$ws.Range("S4").Value = 0.006568746792060225
$ws.Range("T4").Value = 0.006568746792060223
# Row 5
$ws.Range("H5").Value = 63.825936
$ws.Range("I5").Value = 0.6259016025719319
$ws.Range("J5").Value = 0.6259016025719319
$ws.Range("M5").Value = 1.660421
$ws.Range("N5").Value = 4.981262999999999
$ws.Range("O5").Value = 0.03714789785507311
$ws.Range("P5").Value = 0.03714789785507311
$ws.Range("Q5").Value = 35.325974826352
$ws.Range("R5").Value = 317.9337734371679
$ws.Range("S5").Value = 0.02325092879966869
$ws.Range("T5").Value = 0.02325092879966869
# Row 6
$ws.Range("H6").Value = 63.825936
$ws.Range("I6").Value = 0.6259016025719319
$ws.Range("J6").Value = 0.6259016025719319
$ws.Range("O6").Value = 0.5631392661118858
$ws.Range("P6").Value = 0.5631392661118859
$ws.Range("Q6").Value = 535.5200344312906
$ws.Range("R6").Value = 4819.680309881615
$ws.Range("S6").Value = 0.352469769130611
$ws.Range("T6").Value = 0.3524697691306111
# Row 7
$ws.Range("H7").Value = 63.825936
$ws.Range("I7").Value = 0.6259016025719319
$ws.Range("J7").Value = 0.6259016025719319
$ws.Range("M7").Value = 17.866195
$ws.Range("N7").Value = 53.598585
$ws.Range("O7").Value = 0.399712836033041
$ws.Range("P7").Value = 0.399712836033041
$ws.Range("Q7").Value = 380.10887287784
$ws.Range("R7").Value = 3420.97985590056
$ws.Range("S7").Value = 0.2501809046416523
$ws.Range("T7").Value = 0.2501809046416523
# Row 8
$ws.Range("G8").Value = 12.157548
$ws.Range("H8").Value = 36.472644
$ws.Range("I8").Value = 0.3576647325569273
$ws.Range("J8").Value = 0.3576647325569273
$ws.Range("M8").Value = 1.660421
$ws.Range("N8").Value = 4.981262999999999
$ws.Range("O8").Value = 0.03714789785507311
$ws.Range("P8").Value = 0.03714789785507311
$ws.Range("Q8").Value = 20.186648007708
$ws.Range("R8").Value = 181.679832069372
$ws.Range("S8").Value = 0.01328649295138678
$ws.Range("T8").Value = 0.01328649295138678
# Row 9
$ws.Range("G9").Value = 12.157548
$ws.Range("H9").Value = 36.472644
$ws.Range("I9").Value = 0.3576647325569273
$ws.Range("J9").Value = 0.3576647325569273
$ws.Range("O9").Value = 0.5631392661118858
$ws.Range("P9").Value = 0.5631392661118859
$ws.Range("Q9").Value = 306.017158458596
$ws.Range("R9").Value = 2754.154426127364
$ws.Range("S9").Value = 0.201415055006212
$ws.Range("T9").Value = 0.201415055006212
# Row 10
$ws.Range("G10").Value = 12.157548
$ws.Range("H10").Value = 36.472644
$ws.Range("I10").Value = 0.3576647325569273
$ws.Range("J10").Value = 0.3576647325569273
$ws.Range("M10").Value = 17.866195
$ws.Range("N10").Value = 53.598585
$ws.Range("O10").Value = 0.399712836033041
$ws.Range("P10").Value = 0.399712836033041
$ws.Range("Q10").Value = 217.20912328986
$ws.Range("R10").Value = 1954.88210960874
$ws.Range("S10").Value = 0.1429631845993286
$ws.Range("T10").Value = 0.1429631845993286